$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81; existing rows 81..153 shift down to 82..154
$ws.Rows("81:81").Insert()

# Populate the newly inserted row 81 with its data (columns A-T)
$ws.Range("A81").Value = 5
$ws.Range("B81").Value = "Macroferia Regional de Talca"
$ws.Range("C81").Value = "Maule"
$ws.Range("D81").Value = 44484
$ws.Range("E81").Value = 7
$ws.Range("F81").Value = "Fruta"
$ws.Range("G81").Value = 100108
$ws.Range("H81").Value = "Tropicales y subtropicales"
$ws.Range("I81").Value = 100108005
$ws.Range("J81").Value = "Piña"
$ws.Range("K81").Value = "Caramelo"
$ws.Range("L81").Value = "Segunda"
$ws.Range("M81").Value = 150
$ws.Range("N81").Value = 22000
$ws.Range("O81").Value = 22000
$ws.Range("P81").Value = 22000
$ws.Range("Q81").Value = "$/caja 14 unidades"
$ws.Range("R81").Value = "Ecuador"
$ws.Range("S81").Value = 1571
$ws.Range("T81").Value = 14
